$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("주석")
$ws.Select()
$ws.Columns("D").Insert()
